$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column values are written as text (matching original inline-string cells)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.172.28'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '1.637.00'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D5").Value = '216.74'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '20.05'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '1.866.44'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '1.641.05'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = '0.540'
$ws.Range("E15").Value = '  +1.29%  '
$ws.Range("D16").Value = '66.53'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '27.160.80'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = '216.46'
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '6.82'
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").Value = '2.54'
$ws.Range("E22").Value = '  +4.68%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '147.58'
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").Value = '0.118'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").Value = '15.65'
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").Value = '1.302.09'
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '0.551'
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("E42").Value = '  +5.51%  '
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("D44").Value = '1.776.21'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").Value = '91.29'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("D50").Value = '7.62'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("E51").Value = '  -0.78%  '

# Restore default styling on the D column so no stray number-format style remains
$ws.Range("D2:D51").Style = "Normal"
